# Fruta / hortaliza, semanal
# Insert a new weekly record at row 44 (shifting all existing data rows
# 44-145 down to 45-146), then populate the new row with the latest
# observation for Vega Monumental Concepción - Piña.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Shift existing row 44 (and everything below it) down by one row.
$ws.Rows.Item(44).Insert(-4121)

# Populate the newly inserted row 44 with the new weekly observation.
$ws.Cells.Item(44, 1).Value  = 11
$ws.Cells.Item(44, 2).Value  = "Vega Monumental Concepción"
$ws.Cells.Item(44, 3).Value  = "Bíobío"
$ws.Cells.Item(44, 4).Value  = 44614
$ws.Cells.Item(44, 5).Value  = 8
$ws.Cells.Item(44, 6).Value  = "Fruta"
$ws.Cells.Item(44, 7).Value  = 100108
$ws.Cells.Item(44, 8).Value  = "Tropicales y subtropicales"
$ws.Cells.Item(44, 9).Value  = 100108005
$ws.Cells.Item(44, 10).Value = "Piña"
$ws.Cells.Item(44, 11).Value = "Caramelo"
$ws.Cells.Item(44, 12).Value = "Segunda"
$ws.Cells.Item(44, 13).Value = 200
$ws.Cells.Item(44, 14).Value = 16000
$ws.Cells.Item(44, 15).Value = 17000
$ws.Cells.Item(44, 16).Value = 16500
$ws.Cells.Item(44, 17).Value = "`$/caja 14 unidades"
$ws.Cells.Item(44, 18).Value = "Ecuador"
$ws.Cells.Item(44, 19).Value = 1179
$ws.Cells.Item(44, 20).Value = 14

# Make sure the date column keeps its date number format (it should
# already have inherited it from the Insert shift, but set it
# explicitly to be safe).
$ws.Cells.Item(44, 4).NumberFormat = "YYYY-MM-DD HH:MM:SS"
